$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.238.12'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.810.17'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.43'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9991'
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5138'
$ws.Range('E7').Value = '  -2.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3951'
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07812'
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.110'
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.01'
$ws.Range('E11').Value = '  -2.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.355'
$ws.Range('E12').Value = '  -0.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.0000'
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.48'
$ws.Range('E14').Value = '  -1.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.341'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.802.87'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.67'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001079'
$ws.Range('E18').Value = '  -1.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06557'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9985'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.32'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.006'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.291.57'
$ws.Range('E23').Value = '  -0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.13'
$ws.Range('E24').Value = '  -1.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.223'
$ws.Range('E25').Value = '  -1.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.73'
$ws.Range('E26').Value = '  +1.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.456'
$ws.Range('E27').Value = '  +1.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.51'
$ws.Range('E28').Value = '  -1.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.019.01'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.65'
$ws.Range('E30').Value = '  +2.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1098'
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.064'
$ws.Range('E32').Value = '  -1.22%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.649'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.563'
$ws.Range('E34').Value = '  -1.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07171'
$ws.Range('E35').Value = '  -2.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.220'
$ws.Range('E36').Value = '  +5.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02351'
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2176'
$ws.Range('E38').Value = '  -0.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.57'
$ws.Range('E39').Value = '  -5.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.031'
$ws.Range('E40').Value = '  -1.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6180'
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9984'
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.158'
$ws.Range('E43').Value = '  -1.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.27'
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5982'
$ws.Range('E45').Value = '  -1.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.305'
$ws.Range('E46').Value = '  -5.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.740'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.09'
$ws.Range('E48').Value = '  -1.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.208'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.920'
$ws.Range('E50').Value = '  -2.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06796'
$ws.Range('E51').Value = '  -1.43%  '
